# Update "想去人数" (F column) counts to the latest scrape snapshot.
# Mirrors the upstream gh-pages data refresh (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 1655
$ws.Range("F5").Value  = 746
$ws.Range("F6").Value  = 623
$ws.Range("F8").Value  = 325
$ws.Range("F9").Value  = 11
$ws.Range("F11").Value = 1638
$ws.Range("F12").Value = 1413
$ws.Range("F15").Value = 1473
$ws.Range("F20").Value = 83
$ws.Range("F21").Value = 386
$ws.Range("F22").Value = 1109
$ws.Range("F23").Value = 96
$ws.Range("F26").Value = 280
$ws.Range("F28").Value = 261
$ws.Range("F29").Value = 74
$ws.Range("F30").Value = 605
$ws.Range("F31").Value = 634
$ws.Range("F36").Value = 321
$ws.Range("F38").Value = 253
$ws.Range("F39").Value = 621
$ws.Range("F41").Value = 1254
$ws.Range("F45").Value = 319
$ws.Range("F46").Value = 61
$ws.Range("F47").Value = 321
$ws.Range("F48").Value = 65

# Sheet "演出" (Performance)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value  = 37
$ws.Range("F6").Value  = 69
$ws.Range("F11").Value = 692
$ws.Range("F16").Value = 48
$ws.Range("F19").Value = 970
$ws.Range("F20").Value = 32
$ws.Range("F21").Value = 1064
$ws.Range("F23").Value = 655
$ws.Range("F24").Value = 16
$ws.Range("F26").Value = 316
$ws.Range("F30").Value = 20
$ws.Range("F36").Value = 114

# Sheet "本地生活" (Local Life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value  = 366
$ws.Range("F7").Value  = 2379
$ws.Range("F8").Value  = 3607
$ws.Range("F9").Value  = 15
$ws.Range("F11").Value = 72
$ws.Range("F12").Value = 108

# Sheet "全部类型" (All Types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 1655
$ws.Range("F4").Value  = 366
$ws.Range("F6").Value  = 3607
$ws.Range("F7").Value  = 746
$ws.Range("F8").Value  = 72
$ws.Range("F9").Value  = 72
$ws.Range("F10").Value = 623
$ws.Range("F12").Value = 325
$ws.Range("F13").Value = 692
$ws.Range("F14").Value = 1413
$ws.Range("F16").Value = 108
$ws.Range("F17").Value = 108
$ws.Range("F18").Value = 1473
$ws.Range("F22").Value = 1109
$ws.Range("F23").Value = 96
$ws.Range("F26").Value = 48
$ws.Range("F28").Value = 280
$ws.Range("F29").Value = 32
$ws.Range("F31").Value = 261
$ws.Range("F32").Value = 1064
$ws.Range("F34").Value = 605
$ws.Range("F35").Value = 634
$ws.Range("F36").Value = 16
$ws.Range("F39").Value = 316
$ws.Range("F40").Value = 316
$ws.Range("F41").Value = 321
$ws.Range("F43").Value = 253
$ws.Range("F45").Value = 621
$ws.Range("F49").Value = 114
$ws.Range("F50").Value = 319
$ws.Range("F51").Value = 321
